# Homelessness NPA widget and uploader.
# - Rename "Sheet2" -> "Description"
# - Add a new "Status: Improving" row under "Updated" on the Description sheet
# - Freeze panes on the Data sheet (header rows 1-19 / column A) and restore view state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename Sheet2 -> Description
# ---------------------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("Sheet2")
$wsDesc.Name = "Description"

# ---------------------------------------------------------------------------
# 2) Insert a new "Status" / "Improving" row right after the "Updated" row
#    (old row 2 and below shift down to make room)
# ---------------------------------------------------------------------------
$wsDesc.Rows.Item(2).Insert()
$wsDesc.Range("A2").Value = "Status"
$wsDesc.Range("B2").Value = "Improving"

# ---------------------------------------------------------------------------
# 3) Sheet view bookkeeping
# ---------------------------------------------------------------------------

# Data sheet: freeze the header band (rows 1-19, column A) and leave the
# window parked on D3, matching the refreshed dashboard view.
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()
$wsData.Range("A1").Select()
$wsData.Range("B20").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.DisplayGridlines = $true
$wsData.Range("D3").Select()

# Description sheet stays the active tab, cursor parked on the new B2 cell.
$wsDesc.Activate()
$excel.ActiveWindow.DisplayGridlines = $true
$wsDesc.Range("B2").Select()
